# weapon spreadsheets stat change
# increased amour pen for AP and HE
# corrected incorrect shield damage values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prototype Stats")

# --- AP rounds: Shield Damage -0.2 -> -0.5, Armour Penetration 0.2 -> 0.5 ---
$ws.Range("K9").Value = -0.5
$ws.Range("L9").Value = 0.5

$ws.Range("K13").Value = -0.5
$ws.Range("L13").Value = 0.5

$ws.Range("K17").Value = -0.5
$ws.Range("L17").Value = 0.5

$ws.Range("K21").Value = -0.5
$ws.Range("L21").Value = 0.5

# --- HE rounds: Shield Damage -0.2 -> -0.5, Armour Penetration 0.05 -> 0.25 ---
$ws.Range("K10").Value = -0.5
$ws.Range("L10").Value = 0.25

$ws.Range("K14").Value = -0.5
$ws.Range("L14").Value = 0.25

$ws.Range("K18").Value = -0.5
$ws.Range("L18").Value = 0.25

$ws.Range("K22").Value = -0.5
$ws.Range("L22").Value = 0.25

# --- AB rounds: Shield Damage 1.3 -> 1 ---
$ws.Range("K11").Value = 1
$ws.Range("K15").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("K23").Value = 1

# --- CAN rounds: Shield Damage -0.7 -> -1 ---
$ws.Range("K12").Value = -1
$ws.Range("K16").Value = -1
$ws.Range("K20").Value = -1
$ws.Range("K24").Value = -1

# Update the last active selection to match the saved cursor position
$ws.Range("L9").Select()
